$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns per data-row, in order: B, C, D, E, F, I, J, M, N
$cols = @(2, 3, 4, 5, 6, 9, 10, 13, 14)

$data = @(
    @(1.543145160119593, 0.6045642687507211, 0.05575800697532962, 0.09227724228550827, 2.705411538345487, 1.262895537521501, 0.1479057805075676, 0.5724660458667472, 1.799457959988544),
    @(1.44065358144411, 0.5657244472716911, 0.05560136928006898, 0.09237417700506612, 2.681179437664312, 1.258904852297803, 0.1481471940153085, 0.5481573466469314, 1.819625327823779),
    @(1.378553937136871, 0.5422353028170619, 0.0555090124632116, 0.09246258146404607, 2.667926424497182, 1.257205323389073, 0.148358173165196, 0.5335595744811243, 1.832670823019765),
    @(1.353455744877579, 0.5327527696101981, 0.05547234174437143, 0.09250588172473684, 2.662933025250567, 1.256700839107523, 0.1484599138985701, 0.5276930664225716, 1.838153240818659),
    @(1.349300739893238, 0.5311835864494867, 0.05546631106277644, 0.09251351140556352, 2.662128433457397, 1.25662840959648, 0.1484777596385136, 0.5267238951603304, 1.839073624241195),
    @(1.378214613528598, 0.5421070566613366, 0.05550851399297585, 0.092463135956244, 2.667857434524834, 1.25719775911486, 0.1483594814624141, 0.5334801241829936, 1.832744088081434),
    @(1.507633062785544, 0.5910974175980073, 0.05570320845519738, 0.09230467650262497, 2.696718020902466, 1.261363300069327, 0.1479759872093638, 0.5640162006656695, 1.806273922297898),
    @(1.768070890135675, 0.6900508364707321, 0.05611512934375007, 0.09222271856964426, 2.766284719212024, 1.275522736755079, 0.1477226294621801, 0.6265138492158826, 1.759639107451669),
    @(1.963570924416729, 0.7645703173390643, 0.05643595357435061, 0.09230151436985956, 2.825413979438792, 1.289627965807, 0.1478417792980053, 0.674052650673346, 1.728611143821006),
    @(2.05343551081819, 0.7988797615579983, 0.05658582421222214, 0.09236745904887833, 2.854079502260106, 1.296859878005165, 0.1479625768201345, 0.6960376782652844, 1.715202477115533),
    @(2.087600248395177, 0.8119317994283506, 0.05664313792879483, 0.0923967506713268, 2.865190407876781, 1.299716532482833, 0.1480179187247472, 0.7044149169224596, 1.710226916638618),
    @(2.080236239506633, 0.8091181369488822, 0.05663076949965173, 0.09239025020522718, 2.862786066338913, 1.29909603545957, 0.1480055725994944, 0.702608411625917, 1.711293948037373),
    @(2.056243552686396, 0.7999523574407021, 0.05659052822307942, 0.09236978236105031, 2.854988464500821, 1.297092525230795, 0.1479669373343455, 0.6967258350692163, 1.714791087578988),
    @(2.041564950718737, 0.7943458645094097, 0.05656595223130445, 0.09235780753231815, 2.850245585762849, 1.295880720710926, 0.1479445227964931, 0.6931293664170539, 1.716946485492159),
    @(1.957716886246089, 0.7623364372834089, 0.05642623786443757, 0.09229780968311552, 2.823576332608809, 1.289171821980233, 0.1478352268465315, 0.6726231351886724, 1.729501680687214),
    @(1.906517977226486, 0.7428053174048728, 0.05634153052925939, 0.0922687061254841, 2.807669504100289, 1.285265558106985, 0.1477852495938805, 0.6601354893806928, 1.737385071684017),
    @(1.877157244714169, 0.7316100901894629, 0.0562931788501011, 0.09225480018499788, 2.798686559564317, 1.283095474487169, 0.1477627709344986, 0.6529867342390503, 1.741985823284985),
    @(1.86723120349842, 0.727826176107385, 0.05627687142048998, 0.09225057882057719, 2.795673588156632, 1.282373868761525, 0.1477562357353293, 0.6505720860364619, 1.743554954474554),
    @(1.911959121728103, 0.7448804421614454, 0.05635050952968079, 0.09227151103922537, 2.809345595062581, 1.285673442793993, 0.1477899209960398, 0.6614613193624166, 1.736538992232262),
    @(2.063287114089917, 0.8026429396424533, 0.05660233286415384, 0.0923756770916011, 2.857271850155485, 1.297677793196911, 0.1479780247679514, 0.6984522770304267, 1.713761120573537),
    @(2.162975963052133, 0.8407429971727538, 0.05677018220933583, 0.0924689309795923, 2.890086810185664, 1.306212023885024, 0.1481569196576089, 0.722931168545685, 1.699469285952489),
    @(2.109697727189314, 0.8203760790846673, 0.05668029992549783, 0.09241685883084649, 2.872435708274963, 1.301593849354006, 0.1480563124140701, 0.7098384803821034, 1.707042522424743),
    @(1.909498949260865, 0.743942174057338, 0.05634644903996744, 0.092270234134471, 2.808587329555849, 1.285488802465281, 0.1477877895770376, 0.6608618164668059, 1.736921291472257),
    @(1.696893226698592, 0.6629664119413974, 0.0560004882983538, 0.0922204464136307, 2.74606504783317, 1.271045733094724, 0.1477376524003198, 0.6093237394562863, 1.771688500391356)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($excelRow, $cols[$c]).Value = $rowValues[$c]
    }
}
